# Regenerate the handback report: the "a114e3a2-c8c1-4bdc-989e-5580c7f2b62d.md"
# file has now been handed back (in sync) for both zh-cn and de-de, so its
# status / timestamps / stale-error text need to be refreshed across all
# three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$statusSynced = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $statusSynced
$ov.Range("F3").Value = $statusSynced

# --- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $statusSynced
$zh.Range("L3").Value = "2017-01-03 07:38:12"
$zh.Range("R3").Value = ""
$zh.Columns.Item(18).AutoFit()

# --- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $statusSynced
$de.Range("L3").Value = "2017-01-03 07:38:24"
$de.Range("R3").Value = ""
$de.Columns.Item(18).AutoFit()
